$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 12:37:50"
$wsZhCn.Range("H3").Value = "2016-03-21 12:38:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 12:37:53"
$wsDeDe.Range("H3").Value = "2016-03-21 12:38:19"
